# Swap data between row 13 and row 14 for columns:
# A (Id), B (Taxonsorteringsordning), E (TaxonId), F (Artnamn),
# G (Vetenskapligt namn), H (Auktor), Q (Ost), R (Nord)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell13 = $ws.Range($col + "13")
    $cell14 = $ws.Range($col + "14")

    $val13 = $cell13.Value2
    $val14 = $cell14.Value2

    $cell13.Value2 = $val14
    $cell14.Value2 = $val13
}
